$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# Title placeholder: "Price Change Analysis" -> "Branch-wise" + " " + "Analysis" (3 runs)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Branch-wise"
$title.InsertAfter(" ")
$title.InsertAfter("Analysis")

# Subtitle placeholder: "Café Chain Analysis " -> "Plumbing Business Analysis "
$subtitle = $s.Shapes.Item(2).TextFrame.TextRange
$subtitle.Text = "Plumbing Business Analysis "
